$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update evening weight and evening body fat
$ws.Range("C2").Value = 103
$ws.Range("E2").Value = 29.3

# Row 3: morning weight / morning body fat removed;
# evening weight / evening body fat updated to new values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 102.9
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = 29.2

# Row 4: update evening weight and evening body fat
$ws.Range("C4").Value = 102.6
$ws.Range("E4").Value = 28.3

# Row 5: evening weight / evening body fat removed
$ws.Range("C5").ClearContents()
$ws.Range("E5").ClearContents()

# Row 6: update evening weight and evening body fat
$ws.Range("C6").Value = 102.65
$ws.Range("E6").Value = 26.9

# Row 8: evening weight / evening body fat removed
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()

# Drug dosage column (F) is no longer recorded for any day
$ws.Range("F2:F8").ClearContents()
